# Rename the player/guest labels in the first column of the roster table.
#
# The old -> new values form a permutation cycle (e.g. rpinheiro becomes
# guest_jlopes, while spinto becomes rpinheiro, etc.), so a global
# text-based Find/Replace cannot be used safely: an early replacement could
# produce text that a later replacement step would then incorrectly match
# again. Instead we address each table cell positionally (fixed row/column)
# and overwrite its contents directly, which is immune to that problem.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row index (1-based, Word convention) in the roster table -> old/new text
# for column 1 of that row.
$rows = @(
    @{ Row = 2;  Old = "rpinheiro";     New = "guest_jlopes" }
    @{ Row = 3;  Old = "pduarte";       New = "guest_random" }
    @{ Row = 4;  Old = "ggomes";        New = "guest_tymoschuk" }
    @{ Row = 5;  Old = "nsilva";        New = "jsilva" }
    @{ Row = 6;  Old = "guest_random";  New = "nsilva" }
    @{ Row = 7;  Old = "jsilva";        New = "apimenta" }
    @{ Row = 8;  Old = "apimenta";      New = "ggomes" }
    @{ Row = 9;  Old = "guest_cgomes";  New = "guest_ggomes" }
    @{ Row = 10; Old = "guest_jlopes";  New = "guest_jpab" }
    @{ Row = 11; Old = "spinto";        New = "rpinheiro" }
)

foreach ($item in $rows) {
    $cell = $t.Cell($item.Row, 1)
    $cell.Range.Text = $item.New
}
